# Fix lush caves temp
# Adds a missing biome row ("minecraft:lush_caves") with its temperature
# value to the biome_temperature sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("biome_temperature")

# New row directly below the existing last data row (row 37).
$ws.Range("A38").Value = "minecraft:lush_caves"
$ws.Range("B38").Value = 30

# Mirror the author's resulting view state: selection on the newly
# added cell, scrolled so the new row is visible.
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A38").Select()
